$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy original B2-D block (rows 82-121) down to new rows 122-161
# This preserves formatting (fill/border/number-format) and values (group/subject/session/date/time/duration)
$ws.Range("A82:G121").Copy($ws.Range("A122:G161"))

# Step 2: rows 2-41 (B2-B -> B2-A), only Group column changes
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 2).Value = "B2-A"
}

# Step 3: rows 42-81 (B2-C -> B2-B), Group + Date change (date = original value from row r-40)
$ws.Cells.Item(42, 2).Value = "B2-B"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "06/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(43, 2).Value = "B2-B"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "07/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(44, 2).Value = "B2-B"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "08/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(45, 2).Value = "B2-B"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "09/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(46, 2).Value = "B2-B"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "10/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(47, 2).Value = "B2-B"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "13/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(48, 2).Value = "B2-B"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "14/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(49, 2).Value = "B2-B"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "15/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(50, 2).Value = "B2-B"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "16/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(51, 2).Value = "B2-B"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "17/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(52, 2).Value = "B2-B"
$c = $ws.Cells.Item(52, 5)
$c.NumberFormat = "@"
$c.Value = "20/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(53, 2).Value = "B2-B"
$c = $ws.Cells.Item(53, 5)
$c.NumberFormat = "@"
$c.Value = "21/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(54, 2).Value = "B2-B"
$c = $ws.Cells.Item(54, 5)
$c.NumberFormat = "@"
$c.Value = "22/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(55, 2).Value = "B2-B"
$c = $ws.Cells.Item(55, 5)
$c.NumberFormat = "@"
$c.Value = "23/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(56, 2).Value = "B2-B"
$c = $ws.Cells.Item(56, 5)
$c.NumberFormat = "@"
$c.Value = "24/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(57, 2).Value = "B2-B"
$c = $ws.Cells.Item(57, 5)
$c.NumberFormat = "@"
$c.Value = "27/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(58, 2).Value = "B2-B"
$c = $ws.Cells.Item(58, 5)
$c.NumberFormat = "@"
$c.Value = "28/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(59, 2).Value = "B2-B"
$c = $ws.Cells.Item(59, 5)
$c.NumberFormat = "@"
$c.Value = "29/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(60, 2).Value = "B2-B"
$c = $ws.Cells.Item(60, 5)
$c.NumberFormat = "@"
$c.Value = "30/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(61, 2).Value = "B2-B"
$c = $ws.Cells.Item(61, 5)
$c.NumberFormat = "@"
$c.Value = "31/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(62, 2).Value = "B2-B"
$c = $ws.Cells.Item(62, 5)
$c.NumberFormat = "@"
$c.Value = "17/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(63, 2).Value = "B2-B"
$c = $ws.Cells.Item(63, 5)
$c.NumberFormat = "@"
$c.Value = "18/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(64, 2).Value = "B2-B"
$c = $ws.Cells.Item(64, 5)
$c.NumberFormat = "@"
$c.Value = "19/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(65, 2).Value = "B2-B"
$c = $ws.Cells.Item(65, 5)
$c.NumberFormat = "@"
$c.Value = "20/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(66, 2).Value = "B2-B"
$c = $ws.Cells.Item(66, 5)
$c.NumberFormat = "@"
$c.Value = "07/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(67, 2).Value = "B2-B"
$c = $ws.Cells.Item(67, 5)
$c.NumberFormat = "@"
$c.Value = "08/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(68, 2).Value = "B2-B"
$c = $ws.Cells.Item(68, 5)
$c.NumberFormat = "@"
$c.Value = "09/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(69, 2).Value = "B2-B"
$c = $ws.Cells.Item(69, 5)
$c.NumberFormat = "@"
$c.Value = "10/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(70, 2).Value = "B2-B"
$c = $ws.Cells.Item(70, 5)
$c.NumberFormat = "@"
$c.Value = "21/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(71, 2).Value = "B2-B"
$c = $ws.Cells.Item(71, 5)
$c.NumberFormat = "@"
$c.Value = "11/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(72, 2).Value = "B2-B"
$c = $ws.Cells.Item(72, 5)
$c.NumberFormat = "@"
$c.Value = "03/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(73, 2).Value = "B2-B"
$c = $ws.Cells.Item(73, 5)
$c.NumberFormat = "@"
$c.Value = "04/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(74, 2).Value = "B2-B"
$c = $ws.Cells.Item(74, 5)
$c.NumberFormat = "@"
$c.Value = "05/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(75, 2).Value = "B2-B"
$c = $ws.Cells.Item(75, 5)
$c.NumberFormat = "@"
$c.Value = "06/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(76, 2).Value = "B2-B"
$c = $ws.Cells.Item(76, 5)
$c.NumberFormat = "@"
$c.Value = "07/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(77, 2).Value = "B2-B"
$c = $ws.Cells.Item(77, 5)
$c.NumberFormat = "@"
$c.Value = "10/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(78, 2).Value = "B2-B"
$c = $ws.Cells.Item(78, 5)
$c.NumberFormat = "@"
$c.Value = "11/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(79, 2).Value = "B2-B"
$c = $ws.Cells.Item(79, 5)
$c.NumberFormat = "@"
$c.Value = "12/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(80, 2).Value = "B2-B"
$c = $ws.Cells.Item(80, 5)
$c.NumberFormat = "@"
$c.Value = "13/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(81, 2).Value = "B2-B"
$c = $ws.Cells.Item(81, 5)
$c.NumberFormat = "@"
$c.Value = "14/01/2026"
$c.NumberFormat = "dd/mm/yyyy"

# Step 4: rows 82-121 (B2-D -> B2-C), Group + Date change (date = original value from row r-40)
$ws.Cells.Item(82, 2).Value = "B2-C"
$c = $ws.Cells.Item(82, 5)
$c.NumberFormat = "@"
$c.Value = "03/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(83, 2).Value = "B2-C"
$c = $ws.Cells.Item(83, 5)
$c.NumberFormat = "@"
$c.Value = "04/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(84, 2).Value = "B2-C"
$c = $ws.Cells.Item(84, 5)
$c.NumberFormat = "@"
$c.Value = "05/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(85, 2).Value = "B2-C"
$c = $ws.Cells.Item(85, 5)
$c.NumberFormat = "@"
$c.Value = "06/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(86, 2).Value = "B2-C"
$c = $ws.Cells.Item(86, 5)
$c.NumberFormat = "@"
$c.Value = "07/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(87, 2).Value = "B2-C"
$c = $ws.Cells.Item(87, 5)
$c.NumberFormat = "@"
$c.Value = "10/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(88, 2).Value = "B2-C"
$c = $ws.Cells.Item(88, 5)
$c.NumberFormat = "@"
$c.Value = "11/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(89, 2).Value = "B2-C"
$c = $ws.Cells.Item(89, 5)
$c.NumberFormat = "@"
$c.Value = "12/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(90, 2).Value = "B2-C"
$c = $ws.Cells.Item(90, 5)
$c.NumberFormat = "@"
$c.Value = "13/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(91, 2).Value = "B2-C"
$c = $ws.Cells.Item(91, 5)
$c.NumberFormat = "@"
$c.Value = "14/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(92, 2).Value = "B2-C"
$c = $ws.Cells.Item(92, 5)
$c.NumberFormat = "@"
$c.Value = "17/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(93, 2).Value = "B2-C"
$c = $ws.Cells.Item(93, 5)
$c.NumberFormat = "@"
$c.Value = "18/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(94, 2).Value = "B2-C"
$c = $ws.Cells.Item(94, 5)
$c.NumberFormat = "@"
$c.Value = "19/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(95, 2).Value = "B2-C"
$c = $ws.Cells.Item(95, 5)
$c.NumberFormat = "@"
$c.Value = "20/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(96, 2).Value = "B2-C"
$c = $ws.Cells.Item(96, 5)
$c.NumberFormat = "@"
$c.Value = "21/01/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(97, 2).Value = "B2-C"
$c = $ws.Cells.Item(97, 5)
$c.NumberFormat = "@"
$c.Value = "07/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(98, 2).Value = "B2-C"
$c = $ws.Cells.Item(98, 5)
$c.NumberFormat = "@"
$c.Value = "08/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(99, 2).Value = "B2-C"
$c = $ws.Cells.Item(99, 5)
$c.NumberFormat = "@"
$c.Value = "09/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(100, 2).Value = "B2-C"
$c = $ws.Cells.Item(100, 5)
$c.NumberFormat = "@"
$c.Value = "10/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(101, 2).Value = "B2-C"
$c = $ws.Cells.Item(101, 5)
$c.NumberFormat = "@"
$c.Value = "11/02/2026"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(102, 2).Value = "B2-C"
$c = $ws.Cells.Item(102, 5)
$c.NumberFormat = "@"
$c.Value = "20/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(103, 2).Value = "B2-C"
$c = $ws.Cells.Item(103, 5)
$c.NumberFormat = "@"
$c.Value = "21/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(104, 2).Value = "B2-C"
$c = $ws.Cells.Item(104, 5)
$c.NumberFormat = "@"
$c.Value = "22/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(105, 2).Value = "B2-C"
$c = $ws.Cells.Item(105, 5)
$c.NumberFormat = "@"
$c.Value = "23/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(106, 2).Value = "B2-C"
$c = $ws.Cells.Item(106, 5)
$c.NumberFormat = "@"
$c.Value = "27/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(107, 2).Value = "B2-C"
$c = $ws.Cells.Item(107, 5)
$c.NumberFormat = "@"
$c.Value = "28/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(108, 2).Value = "B2-C"
$c = $ws.Cells.Item(108, 5)
$c.NumberFormat = "@"
$c.Value = "29/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(109, 2).Value = "B2-C"
$c = $ws.Cells.Item(109, 5)
$c.NumberFormat = "@"
$c.Value = "30/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(110, 2).Value = "B2-C"
$c = $ws.Cells.Item(110, 5)
$c.NumberFormat = "@"
$c.Value = "24/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(111, 2).Value = "B2-C"
$c = $ws.Cells.Item(111, 5)
$c.NumberFormat = "@"
$c.Value = "31/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(112, 2).Value = "B2-C"
$c = $ws.Cells.Item(112, 5)
$c.NumberFormat = "@"
$c.Value = "06/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(113, 2).Value = "B2-C"
$c = $ws.Cells.Item(113, 5)
$c.NumberFormat = "@"
$c.Value = "07/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(114, 2).Value = "B2-C"
$c = $ws.Cells.Item(114, 5)
$c.NumberFormat = "@"
$c.Value = "08/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(115, 2).Value = "B2-C"
$c = $ws.Cells.Item(115, 5)
$c.NumberFormat = "@"
$c.Value = "09/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(116, 2).Value = "B2-C"
$c = $ws.Cells.Item(116, 5)
$c.NumberFormat = "@"
$c.Value = "10/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(117, 2).Value = "B2-C"
$c = $ws.Cells.Item(117, 5)
$c.NumberFormat = "@"
$c.Value = "13/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(118, 2).Value = "B2-C"
$c = $ws.Cells.Item(118, 5)
$c.NumberFormat = "@"
$c.Value = "14/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(119, 2).Value = "B2-C"
$c = $ws.Cells.Item(119, 5)
$c.NumberFormat = "@"
$c.Value = "15/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(120, 2).Value = "B2-C"
$c = $ws.Cells.Item(120, 5)
$c.NumberFormat = "@"
$c.Value = "16/12/2025"
$c.NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(121, 2).Value = "B2-C"
$c = $ws.Cells.Item(121, 5)
$c.NumberFormat = "@"
$c.Value = "17/12/2025"
$c.NumberFormat = "dd/mm/yyyy"

Write-Output "edit complete"